# Add two new columns I (I0) and J (IF) to the sheet, mirroring the
# formatting already used by the existing header/data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new header cells I1 ("I0") and J1 ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, centered alignment) from the
# existing header cell H1 so the new headers match the rest of the row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows 2-35: new values for columns I and J ---
$values = @(
    @(9,9),
    @(5,6),
    @(5,7),
    @(5,5),
    @(8,8),
    @(7,9),
    @(1,5),
    @(1,5),
    @(6,8),
    @(7,8),
    @(5,5),
    @(6,9),
    @(4,7),
    @(1,4),
    @(1,4),
    @(1,5),
    @(1,6),
    @(1,6),
    @(1,7),
    @(1,5),
    @(4,7),
    @(7,8),
    @(1,6),
    @(1,7),
    @(1,7),
    @(1,7),
    @(1,5),
    @(1,5),
    @(1,5),
    @(4,8),
    @(1,4),
    @(1,4),
    @(1,3),
    @(1,2)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row = $row + 1
}

$excel.CutCopyMode = $false
